# Applies Week 17 data updates (per commit message) across YDS, OFF, DEF, ST, TURNS, PEN sheets
$wb = $excel.ActiveWorkbook

# --- YDS sheet: extend the run-of-play distance strings (R/P by OFF/DEF) ---
$wsYDS = $wb.Worksheets.Item("YDS")
$wsYDS.Range("B2").Value = '5 7 9 2 19 4 6 1 4 7 2 2 2 1 5 6 17 2 4 13 1 1 4 0 1 1 8 2 5 2 5 6 9 1 -2 4 2 2 5 2 5 2 -2 2 2 1 2 2 3 0 5 1 -1 2 4 -1 5 2 6 5 4 4 1 8 5 3 15 4 1 2 1 7 2 6 6 1 3 0 0 0 9 5 6 2 -1 3 11 3 0 7 0 12 14 3 4 3 2 4 11 29 2 1 4 7 -1 4 -2 10 6 4 11 0 11 1 4 3 4 1 -2 2 -2 3 3 11 8 0 0 4 -1 3 3 3 4 1 0 3 3 4 3 1 -1 13 6 9 4 1 14 3 18 -2 5 8 4 10 -2 2 0 9 16 1 -2 1 3 2 4 7 2 4 5 11 4 0 -1 3 13 4 3 9 1 0 2 21 21 2 2 23 9 0 2 4 4 0 8 5 2 1 -2 1 2 8 -1 7 4 11 3 3 3 3 2 3 3 1 1 3 0 12 8 -1 4 11 5 4 2 5 4 5 5 1 4 7 1 1 1 7 11 7 3 4 1 19 4 6 0 3 2 3 3 1 11 4 0 9 0 4 2 10 -1 3 5 4 7 14 3 5 6 5 1 4 1 3 3 2 3 2 9 3 1 8 6 2 1 5 7 2 1 1 5 3 3 5 6 2 4 0 7 9 3 2 2 4 10 6 2 5 5 38 4 48 14 2 10 0 14 4 3 1 1 2 3 1 6 10 5 3 1 19 26 1 7 3 6 6 5 7 -1 2 1 1 6 3 8 2 4 4 1 1 5 3 3 0 0 5 3 3 11 4 5 5 2 4 9 -4 29 4 5 0 3 3 0 3 2 -2 4 5 7 3 0 1 15 8 4 4 3 1 4 0 2 -2 8 5 3 -3 3 1 2 3 2 -1 5 6 6 4 2 -1 1 -2 7 1 4 3 1 1 3 1 0 1 10 3 4 0 12 0 -6 2 -1 2 4 1 3 2 4 2 5 8 1 1 2 11 1 2 5 3 0 2 1 3 1 2 1 3 7 1 13 3 4 3 2 3 0 0 5 -1 24 5 1 6 1 0 6 6 3 0 11 2 1 3 1 2 2 0 6 35 2 -1 0 4 1 4 4 -2 -2 2 4 13 3 1 6 2 2 0 4 5 2 0 0 2 1 4 2 14 4 1 4 1 4 8 3 9 3 3 1 1 10 2 4 8 1 0 2 4 1 9 5 2 7 4 4 2 1 -3 2 3 2 3 7 -3 4 5 8 1 3 7 1 1 -1 2 0 5 -2 3 1 9 0 1 -1 1 0 -1 0 5 1 4 7 1 1 3 7 2 7 4 3 1 3 4 30 2 -5 2 5 0 7 -1 4 2 4 2 4 0 5 7 0 7 5 5 1 7 16 2 2 8 3 -1 4 0 1 6 5 -6 2 10 0 2 2 -1 5 2 9 0 -3 0 8 3 5 -4 7 4 3 2 3 6 2 0 3 4 2 0 8 3 6 3 4 4 3 6 3 -1 0 2 5 4 2 1 0 9 3 4 25 2 4 2 13 1 2 1 2 3 5 15 6 1 5 2 2 4 7 4 2 1 1 2 1 3 3 4 36 -2 14 -4 6 3 -1 2 6 4 2 5 3 6 5 -1 2 6 4 2 4 3 2 8 2 8 -1 5 3 2'
$wsYDS.Range("C2").Value = '3 3 1 1 18 7 9 9 6 8 8 4 4 3 -1 11 2 27 3 4 3 2 9 0 -2 4 -1 4 14 -1 2 0 5 1 11 3 10 1 4 5 1 3 2 4 14 0 5 1 3 4 14 4 3 0 6 1 7 15 22 3 1 2 2 44 13 7 2 1 4 1 3 20 7 2 1 9 3 0 5 2 11 2 0 2 3 2 11 8 2 4 12 25 4 4 1 2 -1 11 1 2 1 6 2 2 5 2 2 8 1 6 4 3 0 4 2 7 7 3 4 2 7 13 2 6 5 3 1 15 4 3 9 2 0 4 11 4 4 14 0 0 11 2 0 6 4 3 4 -4 3 6 2 2 -2 5 18 3 3 0 34 0 2 0 2 -3 5 14 5 2 9 4 2 3 5 20 4 10 94 2 17 3 5 3 5 5 5 2 2 0 8 2 4 4 0 9 7 5 6 0 3 -3 4 9 -2 5 13 1 2 2 12 1 2 12 3 6 8 3 3 11 1 0 0 2 6 2 3 11 4 5 3 0 0 13 2 9 2 6 0 17 3 1 3 1 2 0 4 3 9 -4 4 6 6 1 0 -3 11 11 4 9 1 11 2 10 7 -1 3 6 2 19 3 4 59 3 4 7 2 1 9 7 4 11 3 4 5 3 3 -1 0 8 -2 -3 9 2 2 6 1 0 1 1 0 8 3 7 4 2 4 1 0 -2 3 1 4 19 11 4 4 2 8 3 3 5 3 8 4 12 -2 1 5 2 0 -3 11 2 3 -1 5 3 6 3 2 19 24 -3 8 3 6 4 0 2 80 3 -1 0 3 3 1 4 13 9 11 5 5 12 0 2 5 -2 3 5 9 -2 1 5 2 12 23 5 1 12 2 2 8 6 11 9 -1 3 11 4 1 14 -2 -2 1 3 2 1 14 0 2 2 6 0 3 5 2 1 11 46 2 7 5 1 1 2 9 3 11 7 5 3 7 2 2 4 3 8 -3 2 2 8 52 10 16 28 -3 4 14 -2 0 6 8 9 4 45 0 5 3 1 6 6 -1 1 3 5 6 1 4 11 7 1 3 5 4 9 13 1 3 9 6 7 6 5 0 3 9 6 0 6 -4 1 1 1 9 -1 2 13 -3 21 7 5 10 8 5 -1 5 8 1 6 1 -1 3 4 5 6 14 1 26 1 5 9 0 3 3 5 -3 10 4 9 1 2 1 0 4 13 7 2 2 -1 4 4 12 0 2 14 0 1 6 6 2 3 4 3 1 8 0 17 8 -5 7 9 6 3 0 8 5 41 12 7 7 1 5 1 4 3 8 2 5 -1 4 4 4 12 0 7 0 1 7 4 -1 4 0 -3 3 -3 4 7 5 15 4 0 8 3 1 24 9 3 2 5 1 1 3 2 3 4 3 3 1 2 0 83 2 -1 4 3 13 1 12 -1 10 11 0 1 8 0 14 -2 4 5 11 6 2 5 4 17 1 -1 6 1 4 -1 0 7 -2 6 16 9 2 1 18 2 7 12 0 1 4 -1 8 7 15 5 0 16 2 7 19 2 3 4 6 3 -2 4 6 1 10 4 0 5 6 2 14 3 16 0 1 2 4 4 6 -2 6 0 5 2 0 -2 7 2 2 2 5 -5 0 4 3 0 -1 2 1 11 7 0 16 4 12 5 4 9 -1 8 -1 9 -1 1 2 16 2 2 1 -3 1 19 2 2 -2 8 -1 6 12 2 2 3 1 8 6 11 5 4 1 9 4 2 24 13 2 0 2 5 -3 1 1 1 2 10 8 6 2 1 10 1 -1 6 0 7 6 0 3 8 2 0 7 0 -1 15 11 -1 4 9 3 16 4 1 11 25 4 5 3 6 -3 7 7 3 6 10 0 12 3 9 13 -1 3 -3 5 3 32 7 -3 8 5 9 3 13 -1 3 13 3 0 3 6 3 47 4 3 0 4 2 3 10 3 0 7 12 16 0 4 1 2 -1 8 8 5 1 1 9 3 -1 2 5 -1 10 2 1 2 3 5 15 3 0 17 9 4 13 -1 0 2 8 1 7 1 2 9 2 -2 6 6 2 13 5 1 1 13 16 6 7 -1 0 2 2 6 3 6 1 10 3 6 6 1 8 2 -2 0 37 3 4 5 2 2'
$wsYDS.Range("B3").Value = '12 5 6 11 9 19 7 15 20 8 7 20 10 12 19 31 17 8 6 4 7 13 6 6 13 38 10 6 5 12 9 13 15 9 6 6 15 4 24 11 7 11 23 3 20 17 28 -1 8 20 20 9 6 8 18 20 15 14 3 14 34 -5 16 6 10 8 15 14 26 13 7 6 24 43 4 14 19 0 13 23 15 16 36 -2 4 14 6 3 13 4 17 11 44 9 23 6 19 9 9 22 5 9 11 19 36 4 22 7 18 17 6 16 8 17 6 12 6 11 35 0 2 4 19 5 6 4 22 53 15 7 6 9 5 6 16 5 8 12 8 24 28 14 9 11 8 8 16 3 6 9 30 1 5 -3 15 13 18 8 12 6 8 57 19 19 9 5 7 11 5 4 6 10 7 10 12 77 6 3 7 7 8 12 7 10 2 5 6 9 8 4 5 4 -2 10 15 4 29 1 10 16 13 6 4 44 2 22 3 16 6 25 21 18 9 24 4 14 4 6 7 13 22 9 6 4 4 14 20 4 10 11 14 2 33 12 37 27 -4 10 38 8 19 18 9 40 34 1 64 11 14 19 7 1 28 13 -2 7 1 15 7 28 11 3 9 9 16 5 9 23 18 8 16 8 13 6 5 11 6 3 48 5 8 3 3 10 8 4 15 12 4 3 15 11 9 4 15 6 9 11 11 2 17 26 38 5 12 20 14 9 4 31 7 2 7 -2 14 5 23 6 8 9 5 16 10 10 6 3 7 50 2 12 17 7 26 25 11 6 7 14 13 -3 4 2 29 13 0 22 15 12 20 4 8 42 9 2 14 7 7 12 8 14 6 7 2 11 14 38 20 27 0 4 9 7 7 32 5 29 10 8 7 40 17 25 7 24 12 10 7 2 52 8 7 9 3 29 4 6 4 32 -1 12 13 6 22 -3 13 9 22 17 14 7 2 11 11 8 35 -1 11 4 5 8 30 10 11 1 10 5 5 13 15 11 5 12 5 3 8 7 8 5 8 12 3 6 12 8 10 7 11 7 67 9 16 13 37 7 15 5 -5 27 27 27 3 14 8 6 3 9 10 14 5 8 -4 18 13 9 4 5 5 3 20 11 17 8 1 9 13 8 6 3 14 7 6 12 3 10 10 0 22 5 5 0 -1 18 8 6 5 2 4 5 -3 1 2 10 8 10 5 10 6 8 0 9 26 -1 15 15 3 -5 1 6 7 8 8 10 7 6 45 24 8 8 9 9 23 2 8 31 23 9 9 15 6 7 5 10 6 9 8 15 4 39 15 4 10 3 5 10 7 4 8 10 9 8 12 3 13 10 5 6 2 8 7 6 6 8 20 13 1 20 40 10 11 5 4 -4 4 7 2 6 5 15 5 12 11 2 10 4 14 2 5 7 17 13 30 5 10 7 0 15 8 11 15 15 12 6 9 18 7 11 -2 19 10 7 11 9 11 12 4 17 6 18 7 2 6 8 10 5 6 3 7 22 7 7 1 7 6 12 11 33 2 6 18 3 43 7 4 11 -4 11 4 14 5 36 6 41 18 12 9 6 27 -1 12 23 4 13 6 3 9 -2 2 8 6 15 24 0 8 6 9 4 8 3 20 9 5 6 14'
$wsYDS.Range("C3").Value = '6 6 11 7 15 10 13 2 6 5 4 9 9 2 17 11 7 12 14 7 19 15 1 3 6 -3 4 29 14 4 18 11 1 19 25 11 7 13 5 13 22 5 11 13 6 15 17 10 6 12 6 6 26 7 7 7 15 5 0 14 1 8 14 24 18 29 11 3 8 26 8 5 5 39 23 19 9 25 25 7 6 8 15 3 7 13 3 4 4 13 25 2 13 19 -3 11 15 7 1 11 8 1 51 12 10 5 7 6 8 12 6 18 13 10 7 9 10 8 6 2 13 5 20 22 45 12 7 15 -1 8 11 5 6 23 7 8 9 7 53 2 17 9 8 0 3 8 1 28 36 3 31 15 16 6 15 1 6 -1 9 5 45 6 16 73 -1 7 4 0 5 17 6 12 6 9 5 7 27 16 4 21 -1 5 9 6 10 13 12 20 12 9 9 6 18 11 18 12 6 10 21 2 10 34 17 7 -2 20 9 0 6 13 30 2 42 15 19 9 8 12 16 4 11 14 4 6 10 9 50 10 6 51 17 16 14 12 12 20 2 6 3 7 9 10 4 2 10 9 4 11 16 14 7 7 6 5 18 12 11 21 -5 1 21 5 11 13 10 9 11 13 2 16 0 10 8 17 8 5 5 12 11 14 12 7 9 7 7 10 5 4 6 9 11 12 8 22 16 16 3 2 11 35 18 13 10 7 7 13 7 4 26 11 7 10 15 2 13 9 0 0 24 8 8 5 6 41 5 28 8 15 6 8 17 8 5 6 9 7 5 33 30 -6 9 6 12 7 14 42 20 14 22 7 11 20 4 21 6 8 1 19 4 4 29 12 8 6 11 11 5 5 5 9 52 13 33 15 8 3 -1 27 22 7 6 5 -3 7 15 41 7 9 6 9 5 18 0 13 9 22 18 4 9 3 29 4 11 13 15 14 13 8 7 6 13 6 5 33 5 18 2 7 1 29 12 26 1 13 8 24 12 7 8 10 10 16 25 10 8 30 5 8 11 23 7 5 25 31 9 20 1 16 7 8 5 16 21 2 37 9 5 5 1 19 15 11 6 9 17 23 4 5 7 16 9 3 3 21 20 1 -3 4 10 13 24 3 10 12 9 17 13 7 51 52 6 18 28 10 13 11 1 8 10 6 1 41 5 19 15 16 2 11 12 13 47 9 8 10 5 8 15 14 3 4 15 9 3 11 28 2 11 52 5 68 11 20 6 5 8 17 18 12 3 6 12 10 5 11 18 21 5 6 10 9 1 4 9 11 2 13 14 6 6 7 11 10 10 11 15 1 8 7 12 7 10 8 12 5 11 -1 16 9 5 7 18 -2 4 7 46 11 5 9 9 4 4 8 14 7 11 5 11 6 1 10 22 16 -3 46 3 7 -3 12 4 8 9 14 6 13 2 8 17 24 14 4 17 0 6 1 29 5 8 24 25 16 55 7 13 22 14 1 29 6 1 4 30 11 12 21 9 8 3 6 5 10 11 2 11 7 5 5 7 5 5 11 7 10 16 9 12 25 9 18 12 15 16 14 6 12 17 11 9 17 7 7 6 6 24 13 8 27 3 3 3 12 12 12 1 6 17 27 12 43 6 8 13 45 29'

# --- OFF sheet: Week 17 offensive stat updates ---
$wsOFF = $wb.Worksheets.Item("OFF")
$wsOFF.Range("C2").Value = 310
$wsOFF.Range("F2").Value = 132
$wsOFF.Range("G2").Value = 91
$wsOFF.Range("H2").Value = 8
$wsOFF.Range("J2").Value = 64
$wsOFF.Range("N2").Value = 46
$wsOFF.Range("O2").Value = 36
$wsOFF.Range("B3").Value = 18
$wsOFF.Range("C3").Value = 332
$wsOFF.Range("E3").Value = 63
$wsOFF.Range("F3").Value = 173
$wsOFF.Range("G3").Value = 61
$wsOFF.Range("H3").Value = 60
$wsOFF.Range("I3").Value = 115
$wsOFF.Range("J3").Value = 119
$wsOFF.Range("L3").Value = 508
$wsOFF.Range("M3").Value = 337
$wsOFF.Range("Q3").Value = 933

# --- DEF sheet: Week 17 defensive stat updates ---
$wsDEF = $wb.Worksheets.Item("DEF")
$wsDEF.Range("C2").Value = 451
$wsDEF.Range("E2").Value = 20
$wsDEF.Range("F2").Value = 116
$wsDEF.Range("G2").Value = 125
$wsDEF.Range("I2").Value = 13
$wsDEF.Range("J2").Value = 64
$wsDEF.Range("N2").Value = 35
$wsDEF.Range("O2").Value = 40
$wsDEF.Range("C3").Value = 324
$wsDEF.Range("E3").Value = 56
$wsDEF.Range("F3").Value = 201
$wsDEF.Range("H3").Value = 53
$wsDEF.Range("I3").Value = 104
$wsDEF.Range("J3").Value = 98
$wsDEF.Range("L3").Value = 526
$wsDEF.Range("M3").Value = 358
$wsDEF.Range("Q3").Value = 1038

# --- ST sheet: Week 17 special-teams updates (numeric cells) ---
$wsST = $wb.Worksheets.Item("ST")
$wsST.Range("B2").Value = 142
$wsST.Range("D2").Value = 136
$wsST.Range("F2").Value = 183
$wsST.Range("G2").Value = 165
$wsST.Range("L2").Value = 40

# --- ST sheet: extend the per-kick/punt distance strings ---
$wsST.Range("D3").Value = '40 41 42 52 31 57 48 54 58 43 44 44 51 50 59 41 33 56 56 52 45 67 44 49 45 50 61 46 23 43 52 26 64 52 37 42 38 44 41 40 36 52 54 39 26 50 46 50 57 38 47 49 33 51 42 48 38 39 52 43 46 49 50 56 53 44 47 44 35 44 49 51 43 61 47 0 45 54 49 42 21 55 45 48 44 42 54 43 44 69 30 52 48 63 58 53 59 59 58 45 48 48 36 48 52 50 38 50 39 56 53 32 45 50 50 47 53 60 41 51 54 47 58 35 35 49 55 51 38 53 53 44 39 53 49 41'
$wsST.Range("B4").Value = '58 65 61 54 61 63 62 62 65 62 64 65 56 63 55 64 56 63 65 66 60 65 69 67 68 64 62 50 69 69 58 58 63 60 62 61 62 64 64 70 62 64 58 57 58 60 54 68 55 57 61 68 67 64 67 58 60 60 67 52'
$wsST.Range("D4").Value = '0 0 0 0 0 8 0 0 13 12 0 0 0 6 0 0 0 0 16 4 0 0 0 6 0 -1 0 4 0 0 9 0 0 2 0 0 5 3 0 0 0 9 35 0 0 0 0 21 13 0 7 0 0 0 0 0 0 0 1 0 0 0 7 15 5 0 12 0 0 0 0 9 6 12 -1 0 13 0 6 0 0 0 2 13 0 0 21 0 0 0 0 0 2 0 8 0 0 0 15 0 10 18 0 0 16 0 4 0 0 0 17 0 0 0 7 -1 16 0 0 15 0 0 0 9 22 -2 48 3 0 0 0 9 0 18 11 0'
$wsST.Range("B5").Value = '20 21 20 16 19 28 20 20 40 24 23 18 17 19 18 21 16 39 23 23 16 28 32 23 22 72 16 8 19 64 16 25 16 20 0 18 21 27 31 22 8 22 25 14 26 17 22 29 15 18 25 20 16 17 23 17 27 21 28 7'
$wsST.Range("D5").Value = '0 0 19 10 0 8 1 0 8 0 5 0 13 9 0 0 0 19 0 14 0 0 0 0 0 0 0 13 9 0 7 7 0 0 0 0 8 0 6 11 0 0 0 0 7 0 0 0 0 2 0 8 0 0 0 9 0 0 0 0 6 11 0 0 3 0 0 14 20 17 0 0 0 1 0 0 2 0 0 0 0 10 0 7 0 0 0 0 4 0 2 0 0 16 3 0 0 14 15 0 8 0'
$wsST.Range("B6").Value = '21 22 23 17 26 4 27 29 22 19 20 20 24 21 22 24 0 10 22 17 20 28 15 10 20 19 20 18 11 26 21 23 12 23 22 0 8 19 24 23 18 19 21 27 19 38 17 24 16 16 24 0 26 22 24 23 34 15 19 15 27 44 18 17 22 28 26 22 13 24 11 98 0 23 0 21 18 26 20 13'

# --- TURNS sheet: Week 17 turnover updates ---
$wsTURNS = $wb.Worksheets.Item("TURNS")
$wsTURNS.Range("B3").Value = 14
$wsTURNS.Range("C3").Value = 13
$wsTURNS.Range("D3").Value = 19
$wsTURNS.Range("E3").Value = 11

# --- PEN sheet: Week 17 penalty updates ---
$wsPEN = $wb.Worksheets.Item("PEN")
$wsPEN.Range("B2").Value = 33
$wsPEN.Range("B3").Value = 42
$wsPEN.Range("D4").Value = 13

